$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, A (timestamp serial), B (Notified Production MW), C (Actual Production MW), E (Lookup text)
# Column D (Quarter) is unchanged by this edit.
$rows = @(
    @(2,45918,1,0,"18.09.20251"),
    @(3,45918.01041666666,1,0,"18.09.20252"),
    @(4,45918.02083333334,1,0,"18.09.20253"),
    @(5,45918.03125,1,0,"18.09.20254"),
    @(6,45918.04166666666,1,0,"18.09.20255"),
    @(7,45918.05208333334,1,0,"18.09.20256"),
    @(8,45918.0625,1,0,"18.09.20257"),
    @(9,45918.07291666666,1,0,"18.09.20258"),
    @(10,45918.08333333334,1,0,"18.09.20259"),
    @(11,45918.09375,1,0,"18.09.202510"),
    @(12,45918.10416666666,1,0,"18.09.202511"),
    @(13,45918.11458333334,1,0,"18.09.202512"),
    @(14,45918.125,1,0,"18.09.202513"),
    @(15,45918.13541666666,1,0,"18.09.202514"),
    @(16,45918.14583333334,1,0,"18.09.202515"),
    @(17,45918.15625,1,0,"18.09.202516"),
    @(18,45918.16666666666,6,0,"18.09.202517"),
    @(19,45918.17708333334,6,0,"18.09.202518"),
    @(20,45918.1875,6,0,"18.09.202519"),
    @(21,45918.19791666666,6,0,"18.09.202520"),
    @(22,45918.20833333334,11,0,"18.09.202521"),
    @(23,45918.21875,11,0,"18.09.202522"),
    @(24,45918.22916666666,15,0,"18.09.202523"),
    @(25,45918.23958333334,11,0,"18.09.202524"),
    @(26,45918.25,109,3,"18.09.202525"),
    @(27,45918.26041666666,120,26,"18.09.202526"),
    @(28,45918.27083333334,137,76,"18.09.202527"),
    @(29,45918.28125,160,161,"18.09.202528"),
    @(30,45918.29166666666,590,276,"18.09.202529"),
    @(31,45918.30208333334,635,399,"18.09.202530"),
    @(32,45918.3125,703,524,"18.09.202531"),
    @(33,45918.32291666666,764,620,"18.09.202532"),
    @(34,45918.33333333334,1279,789,"18.09.202533"),
    @(35,45918.34375,1354,923,"18.09.202534"),
    @(36,45918.35416666666,1417,1011,"18.09.202535"),
    @(37,45918.36458333334,1464,1066,"18.09.202536"),
    @(38,45918.375,1835,1170,"18.09.202537"),
    @(39,45918.38541666666,1877,1266,"18.09.202538"),
    @(40,45918.39583333334,1924,1253,"18.09.202539"),
    @(41,45918.40625,1955,1260,"18.09.202540"),
    @(42,45918.41666666666,2121,1284,"18.09.202541"),
    @(43,45918.42708333334,2141,1347,"18.09.202542"),
    @(44,45918.4375,2158,1403,"18.09.202543"),
    @(45,45918.44791666666,2174,1422,"18.09.202544"),
    @(46,45918.45833333334,2207,1423,"18.09.202545"),
    @(47,45918.46875,2214,1384,"18.09.202546"),
    @(48,45918.47916666666,2217,1334,"18.09.202547"),
    @(49,45918.48958333334,2217,1280,"18.09.202548"),
    @(50,45918.5,2198,1335,"18.09.202549"),
    @(51,45918.51041666666,2195,1347,"18.09.202550"),
    @(52,45918.52083333334,2187,1388,"18.09.202551"),
    @(53,45918.53125,2174,1307,"18.09.202552"),
    @(54,45918.54166666666,2005,1270,"18.09.202553"),
    @(55,45918.55208333334,1991,1202,"18.09.202554"),
    @(56,45918.5625,1975,1170,"18.09.202555"),
    @(57,45918.57291666666,1956,1214,"18.09.202556"),
    @(58,45918.58333333334,1775,1186,"18.09.202557"),
    @(59,45918.59375,1753,1117,"18.09.202558"),
    @(60,45918.60416666666,1726,1133,"18.09.202559"),
    @(61,45918.61458333334,1700,1155,"18.09.202560"),
    @(62,45918.625,1448,1077,"18.09.202561"),
    @(63,45918.63541666666,1420,1005,"18.09.202562"),
    @(64,45918.64583333334,1383,987,"18.09.202563"),
    @(65,45918.65625,1342,907,"18.09.202564"),
    @(66,45918.66666666666,927,751,"18.09.202565"),
    @(67,45918.67708333334,883,725,"18.09.202566"),
    @(68,45918.6875,812,617,"18.09.202567"),
    @(69,45918.69791666666,766,532,"18.09.202568"),
    @(70,45918.70833333334,355,367,"18.09.202569"),
    @(71,45918.71875,309,274,"18.09.202570"),
    @(72,45918.72916666666,256,183,"18.09.202571"),
    @(73,45918.73958333334,230,110,"18.09.202572"),
    @(74,45918.75,34,30,"18.09.202573"),
    @(75,45918.76041666666,22,8,"18.09.202574"),
    @(76,45918.77083333334,22,0,"18.09.202575"),
    @(77,45918.78125,21,0,"18.09.202576"),
    @(78,45918.79166666666,10,0,"18.09.202577"),
    @(79,45918.80208333334,10,0,"18.09.202578"),
    @(80,45918.8125,10,0,"18.09.202579"),
    @(81,45918.82291666666,10,0,"18.09.202580"),
    @(82,45918.83333333334,2,0,"18.09.202581"),
    @(83,45918.84375,2,0,"18.09.202582"),
    @(84,45918.85416666666,2,0,"18.09.202583"),
    @(85,45918.86458333334,2,0,"18.09.202584"),
    @(86,45918.875,1,0,"18.09.202585"),
    @(87,45918.88541666666,1,0,"18.09.202586"),
    @(88,45918.89583333334,1,0,"18.09.202587"),
    @(89,45918.90625,1,0,"18.09.202588"),
    @(90,45918.91666666666,1,0,"18.09.202589"),
    @(91,45918.92708333334,1,0,"18.09.202590"),
    @(92,45918.9375,1,0,"18.09.202591"),
    @(93,45918.94791666666,1,0,"18.09.202592"),
    @(94,45918.95833333334,1,0,"18.09.202593"),
    @(95,45918.96875,1,0,"18.09.202594"),
    @(96,45918.97916666666,1,0,"18.09.202595"),
    @(97,45918.98958333334,1,0,"18.09.202596"),
    @(98,45919,1,0,"19.09.20251"),
    @(99,45919.01041666666,1,0,"19.09.20252"),
    @(100,45919.02083333334,1,0,"19.09.20253"),
    @(101,45919.03125,1,0,"19.09.20254"),
    @(102,45919.04166666666,1,0,"19.09.20255"),
    @(103,45919.05208333334,1,0,"19.09.20256"),
    @(104,45919.0625,1,0,"19.09.20257"),
    @(105,45919.07291666666,1,0,"19.09.20258"),
    @(106,45919.08333333334,1,0,"19.09.20259"),
    @(107,45919.09375,1,0,"19.09.202510"),
    @(108,45919.10416666666,1,0,"19.09.202511"),
    @(109,45919.11458333334,1,0,"19.09.202512"),
    @(110,45919.125,1,0,"19.09.202513"),
    @(111,45919.13541666666,1,0,"19.09.202514"),
    @(112,45919.14583333334,1,0,"19.09.202515"),
    @(113,45919.15625,1,0,"19.09.202516"),
    @(114,45919.16666666666,6,0,"19.09.202517"),
    @(115,45919.17708333334,6,0,"19.09.202518"),
    @(116,45919.1875,6,0,"19.09.202519"),
    @(117,45919.19791666666,6,0,"19.09.202520"),
    @(118,45919.20833333334,14,0,"19.09.202521"),
    @(119,45919.21875,14,0,"19.09.202522"),
    @(120,45919.22916666666,14,0,"19.09.202523"),
    @(121,45919.23958333334,14,0,"19.09.202524"),
    @(122,45919.25,95,2,"19.09.202525"),
    @(123,45919.26041666666,106,19,"19.09.202526"),
    @(124,45919.27083333334,124,55,"19.09.202527"),
    @(125,45919.28125,145,126,"19.09.202528"),
    @(126,45919.29166666666,536,232,"19.09.202529"),
    @(127,45919.30208333334,578,348,"19.09.202530"),
    @(128,45919.3125,643,462,"19.09.202531"),
    @(129,45919.32291666666,724,542,"19.09.202532"),
    @(130,45919.33333333334,1235,663,"19.09.202533"),
    @(131,45919.34375,1303,771,"19.09.202534"),
    @(132,45919.35416666666,1352,822,"19.09.202535"),
    @(133,45919.36458333334,1401,901,"19.09.202536"),
    @(134,45919.375,1734,1001,"19.09.202537"),
    @(135,45919.38541666666,1770,1082,"19.09.202538"),
    @(136,45919.39583333334,1808,1121,"19.09.202539"),
    @(137,45919.40625,1844,0,"19.09.202540"),
    @(138,45919.41666666666,2017,0,"19.09.202541"),
    @(139,45919.42708333334,2039,0,"19.09.202542"),
    @(140,45919.4375,2059,0,"19.09.202543"),
    @(141,45919.44791666666,2075,0,"19.09.202544"),
    @(142,45919.45833333334,2157,0,"19.09.202545"),
    @(143,45919.46875,2164,0,"19.09.202546"),
    @(144,45919.47916666666,2170,0,"19.09.202547"),
    @(145,45919.48958333334,2177,0,"19.09.202548"),
    @(146,45919.5,2144,0,"19.09.202549"),
    @(147,45919.51041666666,2148,0,"19.09.202550"),
    @(148,45919.52083333334,2144,0,"19.09.202551"),
    @(149,45919.53125,2136,0,"19.09.202552"),
    @(150,45919.54166666666,1993,0,"19.09.202553"),
    @(151,45919.55208333334,1982,0,"19.09.202554"),
    @(152,45919.5625,1965,0,"19.09.202555"),
    @(153,45919.57291666666,1946,0,"19.09.202556"),
    @(154,45919.58333333334,1786,0,"19.09.202557"),
    @(155,45919.59375,1759,0,"19.09.202558"),
    @(156,45919.60416666666,1729,0,"19.09.202559"),
    @(157,45919.61458333334,1699,0,"19.09.202560"),
    @(158,45919.625,1455,0,"19.09.202561"),
    @(159,45919.63541666666,1415,0,"19.09.202562"),
    @(160,45919.64583333334,1373,0,"19.09.202563"),
    @(161,45919.65625,1332,0,"19.09.202564"),
    @(162,45919.66666666666,936,0,"19.09.202565"),
    @(163,45919.67708333334,886,0,"19.09.202566"),
    @(164,45919.6875,755,0,"19.09.202567"),
    @(165,45919.69791666666,715,0,"19.09.202568"),
    @(166,45919.70833333334,354,0,"19.09.202569"),
    @(167,45919.71875,296,0,"19.09.202570"),
    @(168,45919.72916666666,221,0,"19.09.202571"),
    @(169,45919.73958333334,200,0,"19.09.202572"),
    @(170,45919.75,67,0,"19.09.202573"),
    @(171,45919.76041666666,22,0,"19.09.202574"),
    @(172,45919.77083333334,21,0,"19.09.202575"),
    @(173,45919.78125,20,0,"19.09.202576"),
    @(174,45919.79166666666,11,0,"19.09.202577"),
    @(175,45919.80208333334,11,0,"19.09.202578"),
    @(176,45919.8125,11,0,"19.09.202579"),
    @(177,45919.82291666666,11,0,"19.09.202580"),
    @(178,45919.83333333334,2,0,"19.09.202581"),
    @(179,45919.84375,2,0,"19.09.202582"),
    @(180,45919.85416666666,2,0,"19.09.202583"),
    @(181,45919.86458333334,2,0,"19.09.202584"),
    @(182,45919.875,1,0,"19.09.202585"),
    @(183,45919.88541666666,1,0,"19.09.202586"),
    @(184,45919.89583333334,1,0,"19.09.202587"),
    @(185,45919.90625,1,0,"19.09.202588"),
    @(186,45919.91666666666,1,0,"19.09.202589"),
    @(187,45919.92708333334,1,0,"19.09.202590"),
    @(188,45919.9375,1,0,"19.09.202591"),
    @(189,45919.94791666666,1,0,"19.09.202592"),
    @(190,45919.95833333334,0,0,"19.09.202593"),
    @(191,45919.96875,0,0,"19.09.202594"),
    @(192,45919.97916666666,0,0,"19.09.202595"),
    @(193,45919.98958333334,0,0,"19.09.202596")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
